# Applies the described changes:
#  1. Update the "Status" shared string value from OPTIMAL to TIME_LIMIT (affects all
#     cells in column E that reference that shared string, i.e. E2:E3 after the row
#     deletion below).
#  2. Fix erroneous data in row 2 and row 3 (columns B, C, D).
#  3. Remove rows 4 through 11 (only rows 1-3 remain), shrinking the used range to A1:H3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Status text used throughout column E.
$ws.Range("E2:E3").Value = "TIME_LIMIT"

# Correct row 2 values.
$ws.Range("B2").Value = -636.8209940539369
$ws.Range("C2").Value = 6.8750427571620385
$ws.Range("D2").Value = 3608.498801201

# Correct row 3 values.
$ws.Range("B3").Value = -642.701897941751
$ws.Range("C3").Value = 9.5353207677372
$ws.Range("D3").Value = 3793.139302787

# Remove rows 4-11, leaving only the header and the two corrected data rows.
$ws.Range("A4:H11").EntireRow.Delete()
